# Applies the crypto price/volume refresh described by the commit diff.
# D-column (Price) values are text that often LOOK numeric (e.g. "328.50",
# "1.007"); a plain COM .Value assignment would let Excel auto-coerce those
# to numbers and silently drop significant trailing/embedded zeros, so those
# cells are written as Text ("@") format and then reset back to the default
# "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.410.68'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.79%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.990.27'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -6.43%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '328.50'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.007'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5000'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.89%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4211'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -6.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '52.23'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08884'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.120'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.55%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.22'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -8.51%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.069.48'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.053'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.502'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -7.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '95.87'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -7.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.009'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001105'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06615'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.65'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -9.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.007'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.966'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -6.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '29.459.85'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.70%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.87'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -7.42%  '
$ws.Range('E25').Value = '  -2.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.53'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.58'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -7.68%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.458'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.94%  '
$ws.Range('E29').Value = '  -9.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.61'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.044'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -10.83%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09919'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.85%  '
$ws.Range('E33').Value = '  -13.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.824'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -7.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.788'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.55%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.548'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -10.88%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02456'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.71%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06334'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -8.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.286'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6505'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -9.03%  '
$ws.Range('E41').Value = '  -8.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.2062'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -9.24%  '
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6325'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -9.03%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.42'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -8.70%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.199'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.89%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.283'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.514'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000330'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06993'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.42%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.138'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.74%  '
